$d = $word.ActiveDocument

foreach ($p in $d.Paragraphs) {
    $styleName = $p.Style.NameLocal

    if ($styleName -eq "Heading 1") {
        $p.Format.LeftIndent = 0
    }
    elseif ($styleName -eq "Heading 2") {
        $p.Format.LeftIndent = 10
    }
    elseif ($styleName -eq "Heading 3") {
        $p.Format.LeftIndent = 20
    }
}
